# Weekly data update: insert a new price-report row for Coliflor
# (Terminal Hortofrutícola Agro Chillán) ahead of the existing row 152,
# pushing the subsequent rows (152-180) down to (153-181).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 152 (shifts rows 152:180 -> 153:181).
$ws.Rows(152).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A152").Value = 7
$ws.Range("B152").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C152").Value = "Ñuble"
$ws.Range("D152").Value = 44504
$ws.Range("E152").Value = 16
$ws.Range("F152").Value = 100112008
$ws.Range("G152").Value = "Coliflor"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 300
$ws.Range("K152").Value = 700
$ws.Range("L152").Value = 750
$ws.Range("M152").Value = 725
$ws.Range("N152").Value = "$/unidad"
$ws.Range("O152").Value = "Región del Maule"
$ws.Range("P152").Value = 725
$ws.Range("Q152").Value = 1
$ws.Range("R152").Value = "Hortaliza"
